$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# sua chiet khau cua sale phu va update chien luoc chay tinh luong theo gio
$ws.Range("C3").Value = 7862571.428571429
$ws.Range("C6").Value = 1065000
$ws.Range("C11").Value = 10238000
$ws.Range("C12").Value = 15757190.47619048
$ws.Range("C13").Value = 17939428.57142857
$ws.Range("C15").Value = 53670333.33333333
